$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension to reflect the new data extent (A1:T13)

# Row 2
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Tgfb2"
$ws.Cells.Item(2, 3).Value = "Eng"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 2
$ws.Cells.Item(2, 6).Value = 0.6666666666666666
$ws.Cells.Item(2, 7).Value = 1.325336333333333
$ws.Cells.Item(2, 8).Value = 3.976009
$ws.Cells.Item(2, 9).Value = 0.02918077208126263
$ws.Cells.Item(2, 10).Value = 0.02918077208126263
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 150.1098686666667
$ws.Cells.Item(2, 14).Value = 450.329606
$ws.Cells.Item(2, 15).Value = 0.7276622610660995
$ws.Cells.Item(2, 16).Value = 0.7276622610660997
$ws.Cells.Item(2, 17).Value = 198.9460629358282
$ws.Cells.Item(2, 18).Value = 1790.514566422454
$ws.Cells.Item(2, 19).Value = 0.02123374659230608
$ws.Cells.Item(2, 20).Value = 0.02123374659230608

# Row 3
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Tgfb2"
$ws.Cells.Item(3, 3).Value = "Eng"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 2
$ws.Cells.Item(3, 6).Value = 0.6666666666666666
$ws.Cells.Item(3, 7).Value = 1.325336333333333
$ws.Cells.Item(3, 8).Value = 3.976009
$ws.Cells.Item(3, 9).Value = 0.02918077208126263
$ws.Cells.Item(3, 10).Value = 0.02918077208126263
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 42.32476666666667
$ws.Cells.Item(3, 14).Value = 126.9743
$ws.Cells.Item(3, 15).Value = 0.2051706239258123
$ws.Cells.Item(3, 16).Value = 0.2051706239258124
$ws.Cells.Item(3, 17).Value = 56.09455106318889
$ws.Cells.Item(3, 18).Value = 504.8509595687
$ws.Cells.Item(3, 19).Value = 0.005987037214549579
$ws.Cells.Item(3, 20).Value = 0.00598703721454958

# Row 4
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Tgfb2"
$ws.Cells.Item(4, 3).Value = "Eng"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 2
$ws.Cells.Item(4, 6).Value = 0.6666666666666666
$ws.Cells.Item(4, 7).Value = 1.325336333333333
$ws.Cells.Item(4, 8).Value = 3.976009
$ws.Cells.Item(4, 9).Value = 0.02918077208126263
$ws.Cells.Item(4, 10).Value = 0.02918077208126263
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 4.940565666666667
$ws.Cells.Item(4, 14).Value = 14.821697
$ws.Cells.Item(4, 15).Value = 0.02394954586187395
$ws.Cells.Item(4, 16).Value = 0.02394954586187395
$ws.Cells.Item(4, 17).Value = 6.547911185252556
$ws.Cells.Item(4, 18).Value = 58.931200667273
$ws.Cells.Item(4, 19).Value = 0.0006988662392450902
$ws.Cells.Item(4, 20).Value = 0.0006988662392450902

# Row 5
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Tgfb2"
$ws.Cells.Item(5, 3).Value = "Eng"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = 0.6666666666666666
$ws.Cells.Item(5, 7).Value = 1.325336333333333
$ws.Cells.Item(5, 8).Value = 3.976009
$ws.Cells.Item(5, 9).Value = 0.02918077208126263
$ws.Cells.Item(5, 10).Value = 0.02918077208126263
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 8.915377333333334
$ws.Cells.Item(5, 14).Value = 26.746132
$ws.Cells.Item(5, 15).Value = 0.04321756914621411
$ws.Cells.Item(5, 16).Value = 0.04321756914621412
$ws.Cells.Item(5, 17).Value = 11.81587350524311
$ws.Cells.Item(5, 18).Value = 106.342861547188
$ws.Cells.Item(5, 19).Value = 0.001261122035161882
$ws.Cells.Item(5, 20).Value = 0.001261122035161882

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Tgfb2"
$ws.Cells.Item(6, 3).Value = "Eng"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 20.45485233333333
$ws.Cells.Item(6, 8).Value = 61.364557
$ws.Cells.Item(6, 9).Value = 0.4503674794711605
$ws.Cells.Item(6, 10).Value = 0.4503674794711605
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 150.1098686666667
$ws.Cells.Item(6, 14).Value = 450.329606
$ws.Cells.Item(6, 15).Value = 0.7276622610660995
$ws.Cells.Item(6, 16).Value = 0.7276622610660997
$ws.Cells.Item(6, 17).Value = 3070.475197352727
$ws.Cells.Item(6, 18).Value = 27634.27677617454
$ws.Cells.Item(6, 19).Value = 0.3277154184226248
$ws.Cells.Item(6, 20).Value = 0.3277154184226249

# Row 7
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Tgfb2"
$ws.Cells.Item(7, 3).Value = "Eng"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 20.45485233333333
$ws.Cells.Item(7, 8).Value = 61.364557
$ws.Cells.Item(7, 9).Value = 0.4503674794711605
$ws.Cells.Item(7, 10).Value = 0.4503674794711605
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 42.32476666666667
$ws.Cells.Item(7, 14).Value = 126.9743
$ws.Cells.Item(7, 15).Value = 0.2051706239258123
$ws.Cells.Item(7, 16).Value = 0.2051706239258124
$ws.Cells.Item(7, 17).Value = 865.7468522094555
$ws.Cells.Item(7, 18).Value = 7791.7216698851
$ws.Cells.Item(7, 19).Value = 0.09240217675899347
$ws.Cells.Item(7, 20).Value = 0.09240217675899348

# Row 8
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Tgfb2"
$ws.Cells.Item(8, 3).Value = "Eng"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 20.45485233333333
$ws.Cells.Item(8, 8).Value = 61.364557
$ws.Cells.Item(8, 9).Value = 0.4503674794711605
$ws.Cells.Item(8, 10).Value = 0.4503674794711605
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 4.940565666666667
$ws.Cells.Item(8, 14).Value = 14.821697
$ws.Cells.Item(8, 15).Value = 0.02394954586187395
$ws.Cells.Item(8, 16).Value = 0.02394954586187395
$ws.Cells.Item(8, 17).Value = 101.0585411548032
$ws.Cells.Item(8, 18).Value = 909.5268703932289
$ws.Cells.Item(8, 19).Value = 0.01078609660429113
$ws.Cells.Item(8, 20).Value = 0.01078609660429113

# Row 9
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Tgfb2"
$ws.Cells.Item(9, 3).Value = "Eng"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 20.45485233333333
$ws.Cells.Item(9, 8).Value = 61.364557
$ws.Cells.Item(9, 9).Value = 0.4503674794711605
$ws.Cells.Item(9, 10).Value = 0.4503674794711605
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 8.915377333333334
$ws.Cells.Item(9, 14).Value = 26.746132
$ws.Cells.Item(9, 15).Value = 0.04321756914621411
$ws.Cells.Item(9, 16).Value = 0.04321756914621412
$ws.Cells.Item(9, 17).Value = 182.3627268492804
$ws.Cells.Item(9, 18).Value = 1641.264541643524
$ws.Cells.Item(9, 19).Value = 0.01946378768525104
$ws.Cells.Item(9, 20).Value = 0.01946378768525105

# Row 10
$ws.Cells.Item(10, 1).Value = "sCs"
$ws.Cells.Item(10, 2).Value = "Tgfb2"
$ws.Cells.Item(10, 3).Value = "Eng"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 23.63794933333334
$ws.Cells.Item(10, 8).Value = 70.913848
$ws.Cells.Item(10, 9).Value = 0.5204517484475769
$ws.Cells.Item(10, 10).Value = 0.5204517484475769
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 150.1098686666667
$ws.Cells.Item(10, 14).Value = 450.329606
$ws.Cells.Item(10, 15).Value = 0.7276622610660995
$ws.Cells.Item(10, 16).Value = 0.7276622610660997
$ws.Cells.Item(10, 17).Value = 3548.289469975988
$ws.Cells.Item(10, 18).Value = 31934.60522978389
$ws.Cells.Item(10, 19).Value = 0.3787130960511687
$ws.Cells.Item(10, 20).Value = 0.3787130960511688

# Row 11
$ws.Cells.Item(11, 1).Value = "sCs"
$ws.Cells.Item(11, 2).Value = "Tgfb2"
$ws.Cells.Item(11, 3).Value = "Eng"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 23.63794933333334
$ws.Cells.Item(11, 8).Value = 70.913848
$ws.Cells.Item(11, 9).Value = 0.5204517484475769
$ws.Cells.Item(11, 10).Value = 0.5204517484475769
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 42.32476666666667
$ws.Cells.Item(11, 14).Value = 126.9743
$ws.Cells.Item(11, 15).Value = 0.2051706239258123
$ws.Cells.Item(11, 16).Value = 0.2051706239258124
$ws.Cells.Item(11, 17).Value = 1000.470690011822
$ws.Cells.Item(11, 18).Value = 9004.2362101064
$ws.Cells.Item(11, 19).Value = 0.1067814099522693
$ws.Cells.Item(11, 20).Value = 0.1067814099522693

# Row 12
$ws.Cells.Item(12, 1).Value = "sCs"
$ws.Cells.Item(12, 2).Value = "Tgfb2"
$ws.Cells.Item(12, 3).Value = "Eng"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 23.63794933333334
$ws.Cells.Item(12, 8).Value = 70.913848
$ws.Cells.Item(12, 9).Value = 0.5204517484475769
$ws.Cells.Item(12, 10).Value = 0.5204517484475769
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 4.940565666666667
$ws.Cells.Item(12, 14).Value = 14.821697
$ws.Cells.Item(12, 15).Value = 0.02394954586187395
$ws.Cells.Item(12, 16).Value = 0.02394954586187395
$ws.Cells.Item(12, 17).Value = 116.7848409066729
$ws.Cells.Item(12, 18).Value = 1051.063568160056
$ws.Cells.Item(12, 19).Value = 0.01246458301833773
$ws.Cells.Item(12, 20).Value = 0.01246458301833773

# Row 13
$ws.Cells.Item(13, 1).Value = "sCs"
$ws.Cells.Item(13, 2).Value = "Tgfb2"
$ws.Cells.Item(13, 3).Value = "Eng"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 23.63794933333334
$ws.Cells.Item(13, 8).Value = 70.913848
$ws.Cells.Item(13, 9).Value = 0.5204517484475769
$ws.Cells.Item(13, 10).Value = 0.5204517484475769
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 8.915377333333334
$ws.Cells.Item(13, 14).Value = 26.746132
$ws.Cells.Item(13, 15).Value = 0.04321756914621411
$ws.Cells.Item(13, 16).Value = 0.04321756914621412
$ws.Cells.Item(13, 17).Value = 210.7412376928818
$ws.Cells.Item(13, 18).Value = 1896.671139235936
$ws.Cells.Item(13, 19).Value = 0.02249265942580119
$ws.Cells.Item(13, 20).Value = 0.02249265942580119
